# bugfix num_impressions and policy report outputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# num_impressions (B1) and policy report (B2) value fixes
$ws.Range("B1").Value = 10
$ws.Range("B2").Value = 1234

# Column widths were nudged slightly narrower when the sheet was resaved
$ws.Columns.Item(1).ColumnWidth = 21.1666666666667
$ws.Columns.Item(2).ColumnWidth = 9.16666666666667
$ws.Columns.Item(3).ColumnWidth = 7.33333333333333

# Active cell moved back to B1
$ws.Range("B1").Select()
